# Apply the cryptos.xlsx cell-value updates (Price / Volume(1h) columns)
# for the GitHub Actions data refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.388.15"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.910.22"
$ws.Range("E3").Value = "  -2.45%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "546.55"
$ws.Range("E5").Value = "  -3.94%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.26"
$ws.Range("E6").Value = "  +3.24%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  +1.76%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.902.32"
$ws.Range("E9").Value = "  -2.60%  "
$ws.Range("E10").Value = "  -3.08%  "
$ws.Range("E11").Value = "  -6.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.446"
$ws.Range("E12").Value = "  +1.54%  "
$ws.Range("E13").Value = "  +0.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.70"
$ws.Range("E14").Value = "  +0.76%  "
$ws.Range("E15").Value = "  +0.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.390.78"
$ws.Range("E16").Value = "  -2.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.80"
$ws.Range("E17").Value = "  +5.96%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.909.33"
$ws.Range("E18").Value = "  -2.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "57.379.45"
$ws.Range("E19").Value = "  -4.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "415.92"
$ws.Range("E20").Value = "  -2.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.09"
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.679"
$ws.Range("E22").Value = "  +2.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.93"
$ws.Range("E23").Value = "  -1.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.98"
$ws.Range("E24").Value = "  +0.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "79.55"
$ws.Range("E25").Value = "  +0.66%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("E28").Value = "  -2.60%  "
$ws.Range("E29").Value = "  +2.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.97"
$ws.Range("E30").Value = "  +1.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.12"
$ws.Range("E31").Value = "  +0.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.91"
$ws.Range("E32").Value = "  -4.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0960"
$ws.Range("E33").Value = "  +2.52%  "
$ws.Range("E34").Value = "  +0.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.931"
$ws.Range("E35").Value = "  +0.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.06"
$ws.Range("E36").Value = "  +0.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "47.75"
$ws.Range("E37").Value = "  -4.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.67"
$ws.Range("E38").Value = "  +4.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0672"
$ws.Range("E39").Value = "  +2.46%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.54"
$ws.Range("E40").Value = "  +3.23%  "
$ws.Range("E41").Value = "  -1.52%  "
$ws.Range("E42").Value = "  -3.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "371.57"
$ws.Range("E43").Value = "  -2.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.656.52"
$ws.Range("E44").Value = "  -0.23%  "
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "121.74"
$ws.Range("E46").Value = "  +1.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.237"
$ws.Range("E47").Value = "  +1.08%  "
$ws.Range("E48").Value = "  +1.73%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.96"
$ws.Range("E49").Value = "  -1.83%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.03"
$ws.Range("E50").Value = "  -2.29%  "
$ws.Range("E51").Value = "  -0.11%  "
